# Add a new member (row 4) to the "Data Simpanan" sheet and top up the
# savings figures for the existing member (row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: top up Simpanan Wajib for Januari/Februari + recompute Total ---
$ws.Range("G3").Value = 5000
$ws.Range("H3").Value = 5000
$ws.Range("R3").Value = 65000

# --- Row 4: new member "002" / T.II/WH/0002 ---
# "002" looks numeric, so a plain .Value assignment would get silently
# coerced to the number 2. Build it as a text formula result instead and
# paste back as a value so it lands as literal text "002" (same as the
# existing "001" cell) without touching the cell's number format/style.
$ws.Range("A4").Formula = '=TEXT(2,"000")'
$ws.Range("A4").Copy()
$ws.Range("A4").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B4").Value = "T.II/WH/0002"
$ws.Range("C4").Value = "Yusuf Husain"
$ws.Range("D4").Value = "Tumbang miwan"
$ws.Range("E4").Value = 50000
$ws.Range("F4").Value = 5000
$ws.Range("G4").Value = 5000
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 60000
